$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("N2").Value = 85.83574689470727

# Row 3 updates
$ws.Range("N3").Value = 85.83574689470727

# Row 4 updates
$ws.Range("D4").Value = 4218.3
$ws.Range("F4").Value = 3.87
$ws.Range("I4").Value = 66
$ws.Range("J4").Value = 66
$ws.Range("K4").Value = 64.2
$ws.Range("N4").Value = 85.83574689470727
